# Bypass_Eligibility_Checker.xlsx - "New Test added for Bypass URL checks"
#
# The test URL used by the Bypass_URL sheet is updated to explicitly include
# the HTTPS port (443) in the address shown on the sheet. The hyperlink
# itself (set up previously via Insert Hyperlink) is left pointing at the
# original address, so we only touch the displayed cell value - not the
# hyperlink target - which also means the cell loses the dedicated
# "Hyperlink" cell style it had (since it's no longer generated by the
# Insert Hyperlink flow).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("A2")

# Drop any formatting (incl. the Hyperlink style) before writing the new text
$cell.Clear()
$cell.Value = "https://pins-test.gopro.net:443/selfservice/web/portal/exemption.html"
